$wb = $excel.ActiveWorkbook

$hsqc = $wb.Worksheets.Item("HSQC")
$hmbc = $wb.Worksheets.Item("HMBC")

# --- HMBC sheet: clear the numbering values in column A, rows 19-45 (keep formatting/style) ---
$hmbc.Range("A19:A45").ClearContents()

# --- Switch the active/selected sheet from HSQC to HMBC ---
$hmbc.Activate()
$hmbc.Range("A19:XFD45").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
